$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B column values with refined/re-estimated figures ---
$ws.Cells.Item(2, 2).Value = 2.962002992630005
$ws.Cells.Item(3, 2).Value = 3.005122900009155
$ws.Cells.Item(4, 2).Value = 2.984125137329102
$ws.Cells.Item(5, 2).Value = 3.032219171524048
$ws.Cells.Item(6, 2).Value = 2.814259052276611
$ws.Cells.Item(7, 2).Value = 2.669142007827759
$ws.Cells.Item(8, 2).Value = 2.651050090789795
$ws.Cells.Item(9, 2).Value = 2.636411905288696
$ws.Cells.Item(10, 2).Value = 2.433090209960938
$ws.Cells.Item(12, 2).Value = 2.222220420837402
$ws.Cells.Item(13, 2).Value = 2.270013809204102
$ws.Cells.Item(14, 2).Value = 2.256533861160278
$ws.Cells.Item(15, 2).Value = 2.243213176727295
$ws.Cells.Item(17, 2).Value = 2.453269243240356
$ws.Cells.Item(19, 2).Value = 1.963053822517395
$ws.Cells.Item(20, 2).Value = 2.066593647003174
$ws.Cells.Item(21, 2).Value = 1.881415724754333
$ws.Cells.Item(22, 2).Value = 2.447355270385742
$ws.Cells.Item(23, 2).Value = 2.548131227493286
$ws.Cells.Item(25, 2).Value = 2.574150085449219
$ws.Cells.Item(26, 2).Value = 2.611109495162964
$ws.Cells.Item(27, 2).Value = 2.705683946609497
$ws.Cells.Item(28, 2).Value = 2.633024215698242
$ws.Cells.Item(29, 2).Value = 2.782319068908691
$ws.Cells.Item(31, 2).Value = 2.258062839508057
$ws.Cells.Item(32, 2).Value = 2.244787216186523
$ws.Cells.Item(33, 2).Value = 1.963913083076477
$ws.Cells.Item(34, 2).Value = 1.744187712669373
$ws.Cells.Item(35, 2).Value = 1.472136259078979
$ws.Cells.Item(36, 2).Value = 1.254570722579956
$ws.Cells.Item(37, 2).Value = 1.093175888061523
$ws.Cells.Item(39, 2).Value = 1.865288138389587
$ws.Cells.Item(40, 2).Value = 1.961798191070557
$ws.Cells.Item(41, 2).Value = 2.265710115432739
$ws.Cells.Item(42, 2).Value = 2.352944374084473
$ws.Cells.Item(43, 2).Value = 2.034587860107422
$ws.Cells.Item(44, 2).Value = 1.924052119255066
$ws.Cells.Item(45, 2).Value = 2.114802122116089
$ws.Cells.Item(46, 2).Value = 2.098948955535889
$ws.Cells.Item(47, 2).Value = 2.642067670822144
$ws.Cells.Item(48, 2).Value = 2.930945873260498
$ws.Cells.Item(49, 2).Value = 2.61341381072998
$ws.Cells.Item(50, 2).Value = 2.50513768196106
$ws.Cells.Item(51, 2).Value = 2.181643724441528
$ws.Cells.Item(52, 2).Value = 2.101353168487549
$ws.Cells.Item(53, 2).Value = 2.435364246368408
$ws.Cells.Item(54, 2).Value = 2.388523817062378
$ws.Cells.Item(55, 2).Value = 2.391724824905396
$ws.Cells.Item(56, 2).Value = 2.438619136810303
$ws.Cells.Item(57, 2).Value = 1.762461423873901
$ws.Cells.Item(58, 2).Value = 1.787616729736328
$ws.Cells.Item(59, 2).Value = 1.711974382400513
$ws.Cells.Item(60, 2).Value = 1.479837775230408
$ws.Cells.Item(61, 2).Value = 1.823668837547302
$ws.Cells.Item(62, 2).Value = 1.159204602241516
$ws.Cells.Item(63, 2).Value = 0.9501993060112
$ws.Cells.Item(64, 2).Value = 0.8143872022628784
$ws.Cells.Item(65, 2).Value = 0.6618974208831787
$ws.Cells.Item(66, 2).Value = 1.20978057384491
$ws.Cells.Item(67, 2).Value = 1.583680391311646
$ws.Cells.Item(68, 2).Value = 1.987721681594849
$ws.Cells.Item(69, 2).Value = 2.276660919189453
$ws.Cells.Item(70, 2).Value = 2.248338460922241
$ws.Cells.Item(71, 2).Value = 2.192283153533936
$ws.Cells.Item(73, 2).Value = 1.899697661399841
$ws.Cells.Item(74, 2).Value = 1.889026284217834
$ws.Cells.Item(75, 2).Value = 1.62309741973877
$ws.Cells.Item(76, 2).Value = 1.751937031745911
$ws.Cells.Item(77, 2).Value = 1.740855097770691
$ws.Cells.Item(78, 2).Value = 1.645658254623413
$ws.Cells.Item(79, 2).Value = 1.922862410545349
$ws.Cells.Item(80, 2).Value = 1.740945339202881
$ws.Cells.Item(81, 2).Value = 1.622418165206909
$ws.Cells.Item(82, 2).Value = 1.74537980556488
$ws.Cells.Item(83, 2).Value = 1.777263283729553
$ws.Cells.Item(84, 2).Value = 1.897095680236816
$ws.Cells.Item(85, 2).Value = 2.071506261825562
$ws.Cells.Item(86, 2).Value = 2.142422437667847
$ws.Cells.Item(87, 2).Value = 2.262210845947266
$ws.Cells.Item(88, 2).Value = 2.271121025085449
$ws.Cells.Item(89, 2).Value = 2.197124242782593
$ws.Cells.Item(90, 2).Value = 2.04584789276123
$ws.Cells.Item(91, 2).Value = 1.699123024940491
$ws.Cells.Item(92, 2).Value = 1.595390319824219
$ws.Cells.Item(93, 2).Value = 1.770164847373962
$ws.Cells.Item(94, 2).Value = 2.065795660018921
$ws.Cells.Item(95, 2).Value = 2.265379667282104
$ws.Cells.Item(96, 2).Value = 2.251874208450317
$ws.Cells.Item(97, 2).Value = 2.197809457778931

# --- Append new rows 98-111 with date (col A) and estimate (col B) ---
# Copy formatting (style) from A97 down to the new date cells in column A
$ws.Cells.Item(97, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item(98, 1), $ws.Cells.Item(111, 1)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(98, 1).Value = 43466
$ws.Cells.Item(98, 2).Value = 2.05131721496582
$ws.Cells.Item(99, 1).Value = 43556
$ws.Cells.Item(99, 2).Value = 2.134792566299438
$ws.Cells.Item(100, 1).Value = 43647
$ws.Cells.Item(100, 2).Value = 2.353148937225342
$ws.Cells.Item(101, 1).Value = 43739
$ws.Cells.Item(101, 2).Value = 2.234418869018555
$ws.Cells.Item(102, 1).Value = 43831
$ws.Cells.Item(102, 2).Value = 2.117206573486328
$ws.Cells.Item(103, 1).Value = 43922
$ws.Cells.Item(103, 2).Value = 1.195309042930603
$ws.Cells.Item(104, 1).Value = 44013
$ws.Cells.Item(104, 2).Value = 1.724944114685059
$ws.Cells.Item(105, 1).Value = 44105
$ws.Cells.Item(105, 2).Value = 1.596286773681641
$ws.Cells.Item(106, 1).Value = 44197
$ws.Cells.Item(106, 2).Value = 1.660072803497314
$ws.Cells.Item(107, 1).Value = 44287
$ws.Cells.Item(107, 2).Value = 4.447394847869873
$ws.Cells.Item(108, 1).Value = 44378
$ws.Cells.Item(108, 2).Value = 4.036748886108398
$ws.Cells.Item(109, 1).Value = 44470
$ws.Cells.Item(109, 2).Value = 5.476613521575928
$ws.Cells.Item(110, 1).Value = 44562
$ws.Cells.Item(110, 2).Value = 6.436047077178955
$ws.Cells.Item(111, 1).Value = 44652
$ws.Cells.Item(111, 2).Value = 6.011281490325928
